# "Generate Report for Archive"
# - Update status text "Ready for handoff" -> "In Translation" everywhere it
#   appears (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4).
# - Narrow the now-shorter "Status" columns (Overview E/F, zh-cn C, de-de C)
#   to match the new text's auto-fit width.

$wb = $excel.ActiveWorkbook

# Replace the status text on every sheet (shared-string text, used across
# the Overview rollup columns and each language sheet's Status column).
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Re-fit the Status columns to the new, shorter text. The target OOXML
# column width (13.4101845877511 chars) corresponds to a COM ColumnWidth
# of 12.5 under this engine's character-width quantization.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
